$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column widths: B-G=6, H=28 (like old A), I=27 (like old C) ---
$ws.Columns.Item(2).ColumnWidth = 5.166666666666666
$ws.Columns.Item(3).ColumnWidth = 5.166666666666666
$ws.Columns.Item(4).ColumnWidth = 5.166666666666666
$ws.Columns.Item(5).ColumnWidth = 5.166666666666666
$ws.Columns.Item(6).ColumnWidth = 5.166666666666666
$ws.Columns.Item(7).ColumnWidth = 5.166666666666666
$ws.Columns.Item(8).ColumnWidth = 27.166666666666668
$ws.Columns.Item(9).ColumnWidth = 26.166666666666668

# --- Header row 2: move count/percentage headers, add year headers ---
$ws.Range("B2").Value = 2017
$ws.Range("C2").Value = 2018
$ws.Range("D2").Value = 2019
$ws.Range("E2").Value = 2020
$ws.Range("F2").Value = 2021
$ws.Range("G2").Value = 2022
$ws.Range("H2").Value = "total_count_of_occurrences"
$ws.Range("I2").Value = "percentage_of_occurrences"

# --- Data row 4 ---
$ws.Range("B4").Value = 69
$ws.Range("C4").Value = 111
$ws.Range("D4").Value = 60
$ws.Range("E4").Value = 83
$ws.Range("F4").Value = 62
$ws.Range("G4").Value = 126
$ws.Range("H4").Value = 511
$ws.Range("I4").Value = 32.67263427109975

# --- Data row 5 ---
$ws.Range("B5").Value = 44
$ws.Range("C5").Value = 33
$ws.Range("D5").Value = 21
$ws.Range("E5").Value = 88
$ws.Range("F5").Value = 30
$ws.Range("G5").Value = 78
$ws.Range("H5").Value = 294
$ws.Range("I5").Value = 18.79795396419437

# --- Data row 6 ---
$ws.Range("B6").Value = 92
$ws.Range("C6").Value = 92
$ws.Range("D6").Value = 155
$ws.Range("E6").Value = 71
$ws.Range("F6").Value = 88
$ws.Range("G6").Value = 38
$ws.Range("H6").Value = 536
$ws.Range("I6").Value = 34.27109974424553

# --- Data row 7 ---
$ws.Range("B7").Value = 54
$ws.Range("C7").Value = 23
$ws.Range("D7").Value = 23
$ws.Range("E7").Value = 20
$ws.Range("F7").Value = 83
$ws.Range("G7").Value = 20
$ws.Range("H7").Value = 223
$ws.Range("I7").Value = 14.25831202046036

# --- Header row 9: move count/percentage headers, add year headers ---
$ws.Range("B9").Value = 2017
$ws.Range("C9").Value = 2018
$ws.Range("D9").Value = 2019
$ws.Range("E9").Value = 2020
$ws.Range("F9").Value = 2021
$ws.Range("G9").Value = 2022
$ws.Range("H9").Value = "total_count_of_occurrences"
$ws.Range("I9").Value = "percentage_of_occurrences"

# --- Data row 11 ---
$ws.Range("B11").Value = 56
$ws.Range("C11").Value = 97
$ws.Range("D11").Value = 53
$ws.Range("E11").Value = 89
$ws.Range("F11").Value = 56
$ws.Range("G11").Value = 135
$ws.Range("H11").Value = 486
$ws.Range("I11").Value = 31.07416879795397

# --- Data row 12 ---
$ws.Range("B12").Value = 42
$ws.Range("C12").Value = 45
$ws.Range("D12").Value = 61
$ws.Range("E12").Value = 79
$ws.Range("F12").Value = 19
$ws.Range("G12").Value = 63
$ws.Range("H12").Value = 309
$ws.Range("I12").Value = 19.75703324808184

# --- Data row 13 ---
$ws.Range("B13").Value = 82
$ws.Range("C13").Value = 80
$ws.Range("D13").Value = 82
$ws.Range("E13").Value = 67
$ws.Range("F13").Value = 166
$ws.Range("G13").Value = 54
$ws.Range("H13").Value = 531
$ws.Range("I13").Value = 33.95140664961637

# --- Data row 14 ---
$ws.Range("B14").Value = 79
$ws.Range("C14").Value = 37
$ws.Range("D14").Value = 63
$ws.Range("E14").Value = 27
$ws.Range("F14").Value = 22
$ws.Range("G14").Value = 10
$ws.Range("H14").Value = 238
$ws.Range("I14").Value = 15.21739130434783

# --- Header row 16: move count/percentage headers, add year headers ---
$ws.Range("B16").Value = 2017
$ws.Range("C16").Value = 2018
$ws.Range("D16").Value = 2019
$ws.Range("E16").Value = 2020
$ws.Range("F16").Value = 2021
$ws.Range("G16").Value = 2022
$ws.Range("H16").Value = "total_count_of_occurrences"
$ws.Range("I16").Value = "percentage_of_occurrences"

# --- Data row 18 ---
$ws.Range("B18").Value = 72
$ws.Range("C18").Value = 91
$ws.Range("D18").Value = 85
$ws.Range("E18").Value = 82
$ws.Range("F18").Value = 39
$ws.Range("G18").Value = 116
$ws.Range("H18").Value = 485
$ws.Range("I18").Value = 31.14964675658317

# --- Data row 19 ---
$ws.Range("B19").Value = 21
$ws.Range("C19").Value = 29
$ws.Range("D19").Value = 56
$ws.Range("E19").Value = 81
$ws.Range("F19").Value = 16
$ws.Range("G19").Value = 87
$ws.Range("H19").Value = 290
$ws.Range("I19").Value = 18.62556197816313

# --- Data row 20 ---
$ws.Range("B20").Value = 111
$ws.Range("C20").Value = 94
$ws.Range("D20").Value = 86
$ws.Range("E20").Value = 39
$ws.Range("F20").Value = 142
$ws.Range("G20").Value = 33
$ws.Range("H20").Value = 505
$ws.Range("I20").Value = 32.43416827231856

# --- Data row 21 ---
$ws.Range("B21").Value = 53
$ws.Range("C21").Value = 44
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = 57
$ws.Range("F21").Value = 66
$ws.Range("G21").Value = 26
$ws.Range("H21").Value = 277
$ws.Range("I21").Value = 17.79062299293513

# --- Header row 23: move count/percentage headers, add year headers ---
$ws.Range("B23").Value = 2017
$ws.Range("C23").Value = 2018
$ws.Range("D23").Value = 2019
$ws.Range("E23").Value = 2020
$ws.Range("F23").Value = 2021
$ws.Range("G23").Value = 2022
$ws.Range("H23").Value = "total_count_of_occurrences"
$ws.Range("I23").Value = "percentage_of_occurrences"

# --- Data row 25 ---
$ws.Range("B25").Value = 65
$ws.Range("C25").Value = 68
$ws.Range("D25").Value = 71
$ws.Range("E25").Value = 105
$ws.Range("F25").Value = 121
$ws.Range("G25").Value = 98
$ws.Range("H25").Value = 528
$ws.Range("I25").Value = 33.97683397683397

# --- Data row 26 ---
$ws.Range("B26").Value = 18
$ws.Range("C26").Value = 32
$ws.Range("D26").Value = 36
$ws.Range("E26").Value = 37
$ws.Range("F26").Value = 73
$ws.Range("G26").Value = 54
$ws.Range("H26").Value = 250
$ws.Range("I26").Value = 16.08751608751609

# --- Data row 27 ---
$ws.Range("B27").Value = 123
$ws.Range("C27").Value = 101
$ws.Range("D27").Value = 122
$ws.Range("E27").Value = 108
$ws.Range("F27").Value = 54
$ws.Range("G27").Value = 69
$ws.Range("H27").Value = 577
$ws.Range("I27").Value = 37.12998712998713

# --- Data row 28 ---
$ws.Range("B28").Value = 53
$ws.Range("C28").Value = 58
$ws.Range("D28").Value = 30
$ws.Range("E28").Value = 12
$ws.Range("F28").Value = 10
$ws.Range("G28").Value = 36
$ws.Range("H28").Value = 199
$ws.Range("I28").Value = 12.80566280566281

